# Issue #36 REST server check
# - mark issues 36, 37, 38 as DONE (Status column C)
# - add issue 39: "need a way to kill and restart server from browser" (Priority 2)
# - add issue 40: "pressing hambuger on settings takes you to main" (Type bug)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Issues")

# Mark existing issues 36-38 as DONE in the Status column (C)
$ws.Range("C36").Value = "DONE"
$ws.Range("C37").Value = "DONE"
$ws.Range("C38").Value = "DONE"

# New issue row 39
$ws.Range("A39").Value = 39
$ws.Range("B39").Value = 2
$ws.Range("E39").Value = "need a way to kill and restart server from browser"
$ws.Range("H39").Value = "need a way to kill and restart server from browser"
$ws.Rows.Item(39).RowHeight = 43.5

# New issue row 40
$ws.Range("A40").Value = 40
$ws.Range("H40").Value = "pressing hambuger on settings takes you to main"
$ws.Range("D40").Value = "bug"

# Restore the active selection to E36, as recorded in the saved view state
$ws.Range("E36").Select()
